# Add two new columns, I ("I0") and J ("IF"), to the sheet.
# I is a constant 1 for every data row (except row 35, which is 3).
# J mirrors the existing H ("IP") column value for every data row
# (except row 35, which is 4 instead of H35's value of 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the rest of row 1 (bold, bordered,
# centered) by copying H1's format onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows (2-36) ----------------------------------------------------
$lastRow = 36
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
}

# Row 35 has different values for I and J than the rest of the data.
$ws.Cells.Item(35, 9).Value = 3
$ws.Cells.Item(35, 10).Value = 4
